$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated results (incidence/prevalence code fix) ---
$ws.Range("C2").Value = 132828.73762661038
$ws.Range("D2").Value = 17.600000000000001

$ws.Range("C3").Value = 13951.401747330559

$ws.Range("C4").Value = 9295.8637336234424
$ws.Range("D4").Value = 13

$ws.Range("C5").Value = 109581.47214565639

$ws.Range("C6").Value = 5834.0514826449216

$ws.Range("C7").Value = 24905.501190896346

$ws.Range("C8").Value = 32369.069340874357
$ws.Range("D8").Value = 17

$ws.Range("C9").Value = 37559.356664277228
$ws.Range("D9").Value = 16.8

$ws.Range("C10").Value = 32133.228773388055

$ws.Range("C11").Value = 27.530174529478874
$ws.Range("D11").Value = 18.899999999999999

# --- Column B widened to fit its (now longer) labels ---
$ws.Columns("B").ColumnWidth = 14.8333

# --- Selection moved off the old results range ---
$ws.Range("I6").Select()
